$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "arts" column (R),
# shifting R:AE -> U:AH. Excel's EntireColumn.Insert carries the
# header style (s="1") into the new blank columns automatically.
$ws.Range("R1:T1").EntireColumn.Insert()

# Populate the three newly inserted header cells.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# Populate the three newly inserted data cells for row 2.
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Lower-case the "Unknown" placeholder values in row 2 (D2:J2).
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"
